$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily cumulative death rows appended to the bottom of the table
# (Date, DeathCovid, DeathWithCovid, Total)
$data = @(
    @("4/15/2021", 10970, 2111, 13081),
    @("4/16/2021", 11043, 2121, 13164),
    @("4/17/2021", 11106, 2135, 13241),
    @("4/18/2021", 11172, 2146, 13318),
    @("4/19/2021", 11244, 2150, 13394),
    @("4/20/2021", 11304, 2161, 13465),
    @("4/21/2021", 11357, 2175, 13532),
    @("4/22/2021", 11405, 2185, 13590)
)

$startRow = 182

# Grab the date formatting (style) already used by the existing date column
# and re-use it on the new date cells instead of creating a new number
# format, so the style sheet stays untouched.
$ws.Range("A181").Copy() | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dateCell.Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Scroll the view down to the newly-added rows and select the last new cell
$lastRow = $startRow + $data.Count - 1
$ws.Application.ActiveWindow.ScrollRow = 143
$ws.Range("A$lastRow").Select()
$excel.CutCopyMode = $false
